# Update column F ("dSF") values on Sheet1 to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value  = -2
$ws.Range("F6").Value  = -4
$ws.Range("F7").Value  = -11
$ws.Range("F8").Value  = -7
$ws.Range("F11").Value = -8
$ws.Range("F14").Value = 1
$ws.Range("F18").Value = 5
$ws.Range("F21").Value = 5
